$wb = $excel.ActiveWorkbook

# --- Sheet1 (Neg_Change): update rows 2-13 with new data, then delete old row 14 ---
$ws1 = $wb.Worksheets.Item(1)

# Row 2
$ws1.Range("A2").Value = "SHRIRAMFIN"
$ws1.Range("B2").Value = 938.7
$ws1.Range("C2").Value = 959.95
$ws1.Range("D2").Value = 931.7
$ws1.Range("E2").Value = 957.5
$ws1.Range("F2").Value = 14632455
$ws1.Range("G2").Value = 31894463
$ws1.Range("H2").Value = -0.5412227194419295
$ws1.Range("I2").Value = "SHRIRAMFIN"

# Row 3
$ws1.Range("A3").Value = "GRASIM"
$ws1.Range("B3").Value = 2811
$ws1.Range("C3").Value = 2842.7
$ws1.Range("D3").Value = 2811
$ws1.Range("E3").Value = 2823.1
$ws1.Range("F3").Value = 399009
$ws1.Range("G3").Value = 796031
$ws1.Range("H3").Value = -0.4987519330277339
$ws1.Range("I3").Value = "GRASIM"

# Row 4
$ws1.Range("A4").Value = "TRENT"
$ws1.Range("B4").Value = 4234
$ws1.Range("C4").Value = 4236.1
$ws1.Range("D4").Value = 4176.6
$ws1.Range("E4").Value = 4183
$ws1.Range("F4").Value = 662345
$ws1.Range("G4").Value = 1316448
$ws1.Range("H4").Value = -0.4968696066992392
$ws1.Range("I4").Value = "TRENT"

# Row 5
$ws1.Range("A5").Value = "ICICIGI"
$ws1.Range("B5").Value = 1950
$ws1.Range("C5").Value = 1963.4
$ws1.Range("D5").Value = 1947
$ws1.Range("E5").Value = 1955.5
$ws1.Range("F5").Value = 219269
$ws1.Range("G5").Value = 477977
$ws1.Range("H5").Value = -0.5412561692299002
$ws1.Range("I5").Value = "ICICIGI"

# Row 6
$ws1.Range("A6").Value = "VEDL"
$ws1.Range("B6").Value = 586
$ws1.Range("C6").Value = 590.8
$ws1.Range("D6").Value = 584.6
$ws1.Range("E6").Value = 586.75
$ws1.Range("F6").Value = 7941430
$ws1.Range("G6").Value = 16201067
$ws1.Range("H6").Value = -0.5098205568806055
$ws1.Range("I6").Value = "VEDL"

# Row 7
$ws1.Range("A7").Value = "TVSMOTOR"
$ws1.Range("B7").Value = 3714
$ws1.Range("C7").Value = 3714
$ws1.Range("D7").Value = 3675.5
$ws1.Range("E7").Value = 3685.2
$ws1.Range("F7").Value = 261018
$ws1.Range("G7").Value = 592689
$ws1.Range("H7").Value = -0.559603771961349
$ws1.Range("I7").Value = "TVSMOTOR"

# Row 8
$ws1.Range("A8").Value = "MPHASIS"
$ws1.Range("B8").Value = 2891.8
$ws1.Range("C8").Value = 2945
$ws1.Range("D8").Value = 2863.3
$ws1.Range("E8").Value = 2934
$ws1.Range("F8").Value = 453003
$ws1.Range("G8").Value = 923952
$ws1.Range("H8").Value = -0.5097115434568029
$ws1.Range("I8").Value = "MPHASIS"

# Row 9
$ws1.Range("A9").Value = "ASHOKLEY"
$ws1.Range("B9").Value = 177.08
$ws1.Range("C9").Value = 178.2
$ws1.Range("D9").Value = 175.58
$ws1.Range("E9").Value = 178.1
$ws1.Range("F9").Value = 8542554
$ws1.Range("G9").Value = 17655028
$ws1.Range("H9").Value = -0.5161404445237923
$ws1.Range("I9").Value = "ASHOKLEY"

# Row 10
$ws1.Range("A10").Value = "LUPIN"
$ws1.Range("B10").Value = 2131.8
$ws1.Range("C10").Value = 2131.8
$ws1.Range("D10").Value = 2103.5
$ws1.Range("E10").Value = 2108.5
$ws1.Range("F10").Value = 212827
$ws1.Range("G10").Value = 476091
$ws1.Range("H10").Value = -0.5529699154153302
$ws1.Range("I10").Value = "LUPIN"

# Row 11
$ws1.Range("A11").Value = "BANKINDIA"
$ws1.Range("B11").Value = 143.01
$ws1.Range("C11").Value = 143.55
$ws1.Range("D11").Value = 140.59
$ws1.Range("E11").Value = 140.91
$ws1.Range("F11").Value = 3701448
$ws1.Range("G11").Value = 7693720
$ws1.Range("H11").Value = -0.5189000899434864
$ws1.Range("I11").Value = "BANKINDIA"

# Row 12
$ws1.Range("A12").Value = "DELHIVERY"
$ws1.Range("B12").Value = 410
$ws1.Range("C12").Value = 413.5
$ws1.Range("D12").Value = 407.3
$ws1.Range("E12").Value = 412
$ws1.Range("F12").Value = 1450947
$ws1.Range("G12").Value = 3452045
$ws1.Range("H12").Value = -0.5796847955342413
$ws1.Range("I12").Value = "DELHIVERY"

# Row 13
$ws1.Range("A13").Value = "NBCC"
$ws1.Range("B13").Value = 116.99
$ws1.Range("C13").Value = 117.2
$ws1.Range("D13").Value = 115.62
$ws1.Range("E13").Value = 116.12
$ws1.Range("F13").Value = 6440370
$ws1.Range("G13").Value = 13396669
$ws1.Range("H13").Value = -0.5192558687536432
$ws1.Range("I13").Value = "NBCC"

# Remove the now-obsolete row 14 (AMBER) - shifts dimension from A1:I14 to A1:I13
$ws1.Rows.Item(14).Delete()

# --- Sheet2 (Pos_Change): update rows 2-10 with new data, then add new row 11 ---
$ws2 = $wb.Worksheets.Item(2)

# Row 2
$ws2.Range("A2").Value = "ITC"
$ws2.Range("B2").Value = 402.5
$ws2.Range("C2").Value = 408.9
$ws2.Range("D2").Value = 401.8
$ws2.Range("E2").Value = 408.85
$ws2.Range("F2").Value = 10081997
$ws2.Range("G2").Value = 6847871
$ws2.Range("H2").Value = 0.4722819690966725
$ws2.Range("I2").Value = "ITC"

# Row 3
$ws2.Range("A3").Value = "HINDUNILVR"
$ws2.Range("B3").Value = 2289.3
$ws2.Range("C3").Value = 2310
$ws2.Range("D3").Value = 2280.7
$ws2.Range("E3").Value = 2298.1
$ws2.Range("F3").Value = 1128801
$ws2.Range("G3").Value = 790773
$ws2.Range("H3").Value = 0.4274652776460501
$ws2.Range("I3").Value = "HINDUNILVR"

# Row 4
$ws2.Range("A4").Value = "INDHOTEL"
$ws2.Range("B4").Value = 743.7
$ws2.Range("C4").Value = 743.85
$ws2.Range("D4").Value = 734
$ws2.Range("E4").Value = 738.55
$ws2.Range("F4").Value = 1609752
$ws2.Range("G4").Value = 1140608
$ws2.Range("H4").Value = 0.4113104589832791
$ws2.Range("I4").Value = "INDHOTEL"

# Row 5
$ws2.Range("A5").Value = "RVNL"
$ws2.Range("B5").Value = 334.2
$ws2.Range("C5").Value = 349.5
$ws2.Range("D5").Value = 333
$ws2.Range("E5").Value = 342.6
$ws2.Range("F5").Value = 20825008
$ws2.Range("G5").Value = 13224486
$ws2.Range("H5").Value = 0.5747309952159956
$ws2.Range("I5").Value = "RVNL"

# Row 6
$ws2.Range("A6").Value = "360ONE"
$ws2.Range("B6").Value = 1169.9
$ws2.Range("C6").Value = 1186.6
$ws2.Range("D6").Value = 1158
$ws2.Range("E6").Value = 1177
$ws2.Range("F6").Value = 473322
$ws2.Range("G6").Value = 332048
$ws2.Range("H6").Value = 0.4254625837228352
$ws2.Range("I6").Value = "360ONE"

# Row 7
$ws2.Range("A7").Value = "JSL"
$ws2.Range("B7").Value = 794
$ws2.Range("C7").Value = 803
$ws2.Range("D7").Value = 787.1
$ws2.Range("E7").Value = 798
$ws2.Range("F7").Value = 347649
$ws2.Range("G7").Value = 241787
$ws2.Range("H7").Value = 0.4378316452083859
$ws2.Range("I7").Value = "JSL"

# Row 8
$ws2.Range("A8").Value = "DIXON"
$ws2.Range("B8").Value = 12844
$ws2.Range("C8").Value = 12909
$ws2.Range("D8").Value = 12520
$ws2.Range("E8").Value = 12855
$ws2.Range("F8").Value = 730977
$ws2.Range("G8").Value = 508602
$ws2.Range("H8").Value = 0.4372279306805715
$ws2.Range("I8").Value = "DIXON"

# Row 9
$ws2.Range("A9").Value = "ICICIPRULI"
$ws2.Range("B9").Value = 650
$ws2.Range("C9").Value = 656.2
$ws2.Range("D9").Value = 646
$ws2.Range("E9").Value = 650.9
$ws2.Range("F9").Value = 850478
$ws2.Range("G9").Value = 543562
$ws2.Range("H9").Value = 0.5646384405090864
$ws2.Range("I9").Value = "ICICIPRULI"

# Row 10
$ws2.Range("A10").Value = "HFCL"
$ws2.Range("B10").Value = 64.05
$ws2.Range("C10").Value = 67.14
$ws2.Range("D10").Value = 63.73
$ws2.Range("E10").Value = 66.2
$ws2.Range("F10").Value = 20754927
$ws2.Range("G10").Value = 13302618
$ws2.Range("H10").Value = 0.5602137113160732
$ws2.Range("I10").Value = "HFCL"

# Row 11
$ws2.Range("A11").Value = "TITAGARH"
$ws2.Range("B11").Value = 820.35
$ws2.Range("C11").Value = 859.4
$ws2.Range("D11").Value = 813.15
$ws2.Range("E11").Value = 835
$ws2.Range("F11").Value = 3146796
$ws2.Range("G11").Value = 2195944
$ws2.Range("H11").Value = 0.4330037560156361
$ws2.Range("I11").Value = "TITAGARH"

